$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I2").Value = "b"
$ws.Range("J2").Value = "Acknowledge (Backchannel)"
$ws.Range("I9").Value = "sv"
$ws.Range("J9").Value = "Statement-opinion"
$ws.Range("I22").Value = "sd"
$ws.Range("J22").Value = "Statement-non-opinion"
$ws.Range("I26").Value = "b"
$ws.Range("J26").Value = "Acknowledge (Backchannel)"
$ws.Range("I28").Value = "b"
$ws.Range("J28").Value = "Acknowledge (Backchannel)"
$ws.Range("I31").Value = "b"
$ws.Range("J31").Value = "Acknowledge (Backchannel)"
$ws.Range("I49").Value = "%"
$ws.Range("J49").Value = "Uninterpretable"
$ws.Range("I52").Value = "sd"
$ws.Range("J52").Value = "Statement-non-opinion"
$ws.Range("I53").Value = "sd"
$ws.Range("J53").Value = "Statement-non-opinion"
$ws.Range("I56").Value = "sv"
$ws.Range("J56").Value = "Statement-opinion"
$ws.Range("I57").Value = "sd"
$ws.Range("J57").Value = "Statement-non-opinion"
$ws.Range("I66").Value = "%"
$ws.Range("J66").Value = "Uninterpretable"
$ws.Range("I76").Value = "sd"
$ws.Range("J76").Value = "Statement-non-opinion"
$ws.Range("I83").Value = "sd"
$ws.Range("J83").Value = "Statement-non-opinion"
$ws.Range("I86").Value = "sd"
$ws.Range("J86").Value = "Statement-non-opinion"
$ws.Range("I102").Value = "aa"
$ws.Range("J102").Value = "Agree/Accept"
$ws.Range("I117").Value = "sd"
$ws.Range("J117").Value = "Statement-non-opinion"
$ws.Range("I130").Value = "sd"
$ws.Range("J130").Value = "Statement-non-opinion"
$ws.Range("I150").Value = "b"
$ws.Range("J150").Value = "Acknowledge (Backchannel)"
$ws.Range("I152").Value = "b"
$ws.Range("J152").Value = "Acknowledge (Backchannel)"
$ws.Range("I153").Value = "b"
$ws.Range("J153").Value = "Acknowledge (Backchannel)"
$ws.Range("I165").Value = "aa"
$ws.Range("J165").Value = "Agree/Accept"
$ws.Range("I177").Value = "sd"
$ws.Range("J177").Value = "Statement-non-opinion"
$ws.Range("I179").Value = "aa"
$ws.Range("J179").Value = "Agree/Accept"
$ws.Range("I187").Value = "aa"
$ws.Range("J187").Value = "Agree/Accept"
$ws.Range("I188").Value = "aa"
$ws.Range("J188").Value = "Agree/Accept"
$ws.Range("I203").Value = "%"
$ws.Range("J203").Value = "Uninterpretable"
$ws.Range("I209").Value = "sd"
$ws.Range("J209").Value = "Statement-non-opinion"
$ws.Range("I216").Value = "sd"
$ws.Range("J216").Value = "Statement-non-opinion"
